# Figure 5 poster slide: nudge the bottom-right "Group 106" figure group
# (the group holding the 4th mini scatter-plot panel + its callouts) to the
# left, nudging it closer to the "Group 105" panel above it so the layout
# reads as a finished 2x2 poster block.
#
# Only the group's own <a:off x="..."/> changes (EMU 5319381 -> 5177571,
# i.e. -141810 EMU sideways); its y offset, ext, chOff and chExt, and every
# child shape inside the group, are left untouched.

$EMU_PER_POINT = 12700

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$grp = $s.Shapes.Item("Group 106")

$newLeftEmu = 5177571
$grp.Left = $newLeftEmu / $EMU_PER_POINT
